$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 106, pushing existing rows 106:149 down to 108:151
$ws.Rows("106:107").Insert()

# Row 106 - new "Primera" entry
$ws.Cells.Item(106, 1).Value = 2
$ws.Cells.Item(106, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(106, 3).Value = "Coquimbo"
$ws.Cells.Item(106, 4).Value = 44636
$ws.Cells.Item(106, 5).Value = 4
$ws.Cells.Item(106, 6).Value = 100112043
$ws.Cells.Item(106, 7).Value = "Pepino ensalada"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 200
$ws.Cells.Item(106, 11).Value = 16000
$ws.Cells.Item(106, 12).Value = 17000
$ws.Cells.Item(106, 13).Value = 16500
$ws.Cells.Item(106, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(106, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(106, 16).Value = 236
$ws.Cells.Item(106, 17).Value = 70
$ws.Cells.Item(106, 18).Value = "Hortaliza"

# Row 107 - new "Segunda" entry
$ws.Cells.Item(107, 1).Value = 2
$ws.Cells.Item(107, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(107, 3).Value = "Coquimbo"
$ws.Cells.Item(107, 4).Value = 44636
$ws.Cells.Item(107, 5).Value = 4
$ws.Cells.Item(107, 6).Value = 100112043
$ws.Cells.Item(107, 7).Value = "Pepino ensalada"
$ws.Cells.Item(107, 8).Value = "Sin especificar"
$ws.Cells.Item(107, 9).Value = "Segunda"
$ws.Cells.Item(107, 10).Value = 240
$ws.Cells.Item(107, 11).Value = 14000
$ws.Cells.Item(107, 12).Value = 15000
$ws.Cells.Item(107, 13).Value = 14500
$ws.Cells.Item(107, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(107, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(107, 16).Value = 145
$ws.Cells.Item(107, 17).Value = 100
$ws.Cells.Item(107, 18).Value = "Hortaliza"
